$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 181
$ws.Range("I6").Value = 181
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 543
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = -431

$ws.Range("H8").Value = 29.75
$ws.Range("I8").Value = 29.75
$ws.Range("K8").Value = 89.25
$ws.Range("M8").Value = 49.75

$ws.Range("H9").Value = 308.66666
$ws.Range("J9").Value = 361.5
$ws.Range("L9").Value = 361.5
$ws.Range("N9").Value = -699.5

$ws.Range("H28").Value = 1264.2142
$ws.Range("I28").Value = 927.63635
$ws.Range("K28").Value = 927.63635
$ws.Range("M28").Value = -442.63635

$ws.Range("H33").Value = 153.23077
$ws.Range("I33").Value = 152
$ws.Range("K33").Value = 152
$ws.Range("M33").Value = 77

$ws.Range("H51").Value = 3637.5
$ws.Range("J51").Value = 3637.5
$ws.Range("L51").Value = 3637.5
$ws.Range("N51").Value = -4605.5

$ws.Range("H101").Value = 754.1429000000001
$ws.Range("I101").Value = 796.5
$ws.Range("K101").Value = 2389.5
$ws.Range("M101").Value = -767.5

$ws.Range("H111").Value = 1804.1666
$ws.Range("I111").Value = 1968.75
$ws.Range("J111").Value = 1475
$ws.Range("K111").Value = 5906.25
$ws.Range("L111").Value = 4425
$ws.Range("M111").Value = -2839.25
$ws.Range("N111").Value = -10559

$ws.Range("H113").Value = 2916.6667
$ws.Range("J113").Value = 3750
$ws.Range("L113").Value = 3750
$ws.Range("N113").Value = -10258

$ws.Range("H118").Value = 299
$ws.Range("I118").Value = 299.25
$ws.Range("K118").Value = 897.75
$ws.Range("M118").Value = 759.25

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = ""
$ws.Range("N136").Value = 0

$ws.Range("H138").Value = 4437.394
$ws.Range("J138").Value = 6245.1377
$ws.Range("L138").Value = 18735.4131
$ws.Range("N138").Value = -29015.4131

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H35").Value = 1390.25
$ws.Range("I35").Value = 1423.8572
$ws.Range("K35").Value = 1423.8572
$ws.Range("M35").Value = -1017.8572

$ws.Range("H94").Value = 90000
$ws.Range("I94").Value = 90000
$ws.Range("K94").Value = 90000
$ws.Range("M94").Value = -89099

$ws.Range("H122").Value = 2045.0834
$ws.Range("I122").Value = 1947
$ws.Range("J122").Value = 2241.25
$ws.Range("K122").Value = 5841
$ws.Range("L122").Value = 6723.75
$ws.Range("M122").Value = -3391
$ws.Range("N122").Value = -11623.75

$ws.Range("H132").Value = 945.3200000000001
$ws.Range("I132").Value = 919.6829
$ws.Range("K132").Value = 2759.0487
$ws.Range("M132").Value = -229.0487000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1749
$ws.Range("I107").Value = 1749
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1749
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = ""
$ws.Range("N107").Value = 171

$ws.Range("H129").Value = 74113.336
$ws.Range("J129").Value = 74113.336
$ws.Range("L129").Value = 74113.336
$ws.Range("N129").Value = -84113.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 121.55556
$ws.Range("I7").Value = 67.61539
$ws.Range("K7").Value = 67.61539
$ws.Range("M7").Value = 45.38461

$ws.Range("H16").Value = 4507.3335
$ws.Range("I16").Value = 4507.3335
$ws.Range("K16").Value = 4507.3335
$ws.Range("M16").Value = -4220.3335

$ws.Range("H107").Value = 812.4706
$ws.Range("I107").Value = 531.5833
$ws.Range("J107").Value = 1486.6
$ws.Range("K107").Value = 531.5833
$ws.Range("L107").Value = 1486.6
$ws.Range("M107").Value = 1388.4167
$ws.Range("N107").Value = -5326.6

$ws.Range("H113").Value = 4507.3335
$ws.Range("I113").Value = 4507.3335
$ws.Range("K113").Value = 4507.3335
$ws.Range("M113").Value = -2337.3335

$ws.Range("H122").Value = 2642.0527
$ws.Range("J122").Value = 2416.3333
$ws.Range("L122").Value = 7248.999899999999
$ws.Range("N122").Value = -12148.9999

$ws.Range("H132").Value = 3499
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1217.1666
$ws.Range("J5").Value = 1389.4286
$ws.Range("L5").Value = 4168.2858
$ws.Range("N5").Value = -4392.2858

$ws.Range("H50").Value = 1486.8572
$ws.Range("I50").Value = 10
$ws.Range("J50").Value = 1733
$ws.Range("K50").Value = 30
$ws.Range("L50").Value = 5199
$ws.Range("M50").Value = 451
$ws.Range("N50").Value = -6161

$ws.Range("H53").Value = 1486.8572
$ws.Range("I53").Value = 10
$ws.Range("J53").Value = 1733
$ws.Range("K53").Value = 30
$ws.Range("L53").Value = 5199
$ws.Range("M53").Value = 451
$ws.Range("N53").Value = -6161

$ws.Range("H94").Value = 5000
$ws.Range("J94").Value = 5000
$ws.Range("L94").Value = 15000
$ws.Range("N94").Value = -16352

$ws.Range("H131").Value = 1462.2858
$ws.Range("I131").Value = 720
$ws.Range("J131").Value = 1759.2
$ws.Range("K131").Value = 2160
$ws.Range("L131").Value = 5277.6
$ws.Range("M131").Value = 2880
$ws.Range("N131").Value = -15357.6

$ws.Range("H135").Value = 1217.1666
$ws.Range("J135").Value = 1389.4286
$ws.Range("L135").Value = 12504.8574
$ws.Range("N135").Value = -17574.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 1980.2
$ws.Range("I41").Value = 2633.6667
$ws.Range("K41").Value = 2633.6667
$ws.Range("M41").Value = -2278.6667

$ws.Range("H80").Value = 51500
$ws.Range("I80").Value = 8000
$ws.Range("K80").Value = 8000
$ws.Range("M80").Value = -7002

$ws.Range("H83").Value = 51500
$ws.Range("I83").Value = 8000
$ws.Range("K83").Value = 40000
$ws.Range("M83").Value = -35008

$ws.Range("H97").Value = 964.2727
$ws.Range("I97").Value = 913.5
$ws.Range("J97").Value = 1099.6666
$ws.Range("K97").Value = 913.5
$ws.Range("L97").Value = 1099.6666
$ws.Range("M97").Value = -417.5
$ws.Range("N97").Value = -2091.6666

$ws.Range("H122").Value = 145296.14
$ws.Range("J122").Value = 501249.5
$ws.Range("L122").Value = 1503748.5
$ws.Range("N122").Value = -1508648.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = 0

$ws.Range("H30").Value = 1866
$ws.Range("I30").Value = 1866
$ws.Range("K30").Value = 1866
$ws.Range("M30").Value = -1758

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = ""
$ws.Range("M126").Value = ""
$ws.Range("N126").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 32023.334
$ws.Range("J54").Value = 23000
$ws.Range("L54").Value = 23000
$ws.Range("N54").Value = -24040

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").Value = ""

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").Value = ""

$ws.Range("H96").Value = 1000
$ws.Range("J96").Value = 1000
$ws.Range("L96").Value = 1000
$ws.Range("N96").Value = -3746
